$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prime the style for the new H column cells on existing rows 110/111 ---
# Copy the existing "s=3" style (used throughout column A-G of this verb table)
# from G111 onto the new H110/H111 cells before writing their values, so the
# new cells land on the same cellXfs entry instead of minting a new one.
$ws.Range("G111").Copy()
$ws.Range("H110:H111").PasteSpecial(-4122)

# --- Row 112 ("怒る" - to get angry) ---
$ws.Range("A111:G111").Copy()
$ws.Range("A112:H112").PasteSpecial(-4122)

# --- Row 113 ("生きる" - to live) ---
$ws.Range("A112:H112").Copy()
$ws.Range("A113:H113").PasteSpecial(-4122)

# --- Values, written in the same left-to-right / top-to-bottom order as the
#     source edit so newly minted shared-string indices line up with the
#     target workbook's <sst> table ---
$ws.Range("A112").Value = "怒る"
$ws.Range("B112").Value = "怒って"
$ws.Range("C112").Value = "怒った"
$ws.Range("D112").Value = "怒らない"
$ws.Range("E112").Value = "怒ります"
$ws.Range("F112").Value = "怒ろう"
$ws.Range("G112").Value = "怒れる"
$ws.Range("H112").Value = "怒られる"

$ws.Range("H111").Value = "巻き込まれる"
$ws.Range("H110").Value = "行われる"

$ws.Range("A113").Value = "生きる"
$ws.Range("B113").Value = "生きて"
$ws.Range("C113").Value = "生きた"
$ws.Range("D113").Value = "生きない"
$ws.Range("E113").Value = "生きます"
$ws.Range("F113").Value = "生きよう"
$ws.Range("G113").Value = "生きられる"
$ws.Range("H113").Value = "生きられる"

# --- Selection / scroll position matching the saved view ---
$ws.Range("G108").Select()

Write-Output "done"
